# "Added last minute updates"
#
# 1. The placeholder bookmark ID in the first paragraph is renamed from the
#    generic AFFARS topic id to the SMC PGI specific id.
# 2. The stray trailing-space run that used to follow that text is removed.
# 3. The first paragraph picks up the same paragraph border + left indent
#    already used by the body paragraph further down the document.

$d = $word.ActiveDocument

$oldId = "**ID__AFFARS_pgi_5301_topic_42__ID**"
$newId = "**ID__AFFARS_SMC_PGI_5301_304__ID**"

# 1. Swap the placeholder id text (wdReplaceOne = 1 -> first match only).
$d.Content.Find.Execute($oldId, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newId, 1) | Out-Null

# 2. The paragraph used to hold the id text plus a second run containing a
# single trailing space; drop that now-orphaned space run.
$p1 = $d.Paragraphs(1)
$paraRange = $p1.Range
$tailStart = $paraRange.Start + $newId.Length
$tailEnd = $paraRange.End - 1
if ($tailEnd -gt $tailStart) {
    $trailingSpace = $d.Range($tailStart, $tailEnd)
    $trailingSpace.Delete()
}

# 3. Give the paragraph a border (5pt space on every side) and widen its
# left indent from 120 twips (6pt) to 225 twips (11.25pt), matching the
# formatting already used elsewhere in the document.
$p1.Range.ParagraphFormat.LeftIndent = 11.25
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

Write-Output ("Paragraph 1 now reads: " + $p1.Range.Text)
